# Update the LR-pairs data with new TPM-derived values.
# The "Target cluster" column (D) changes from "Inflammatory-Mac" to
# "Resolving-Mac" for rows 2-6, and the dependent expression / specificity
# metrics are recalculated accordingly (columns E-T where affected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 0.5023523333333334
$ws.Range("H2").Value = 1.507057
$ws.Range("I2").Value = 0.06515888850144765
$ws.Range("J2").Value = 0.06515888850144765
$ws.Range("M2").Value = 0.07218766666666666
$ws.Range("N2").Value = 0.216563
$ws.Range("Q2").Value = 0.03626364278788889
$ws.Range("R2").Value = 0.326372785091
$ws.Range("S2").Value = 0.06515888850144765
$ws.Range("T2").Value = 0.06515888850144765

# Row 3
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("I3").Value = 0.07275905893716338
$ws.Range("J3").Value = 0.07275905893716339
$ws.Range("M3").Value = 0.07218766666666666
$ws.Range("N3").Value = 0.216563
$ws.Range("Q3").Value = 0.04049345505366667
$ws.Range("R3").Value = 0.364441095483
$ws.Range("S3").Value = 0.07275905893716338
$ws.Range("T3").Value = 0.07275905893716339

# Row 4
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 2.845667666666667
$ws.Range("H4").Value = 8.537003
$ws.Range("I4").Value = 0.3691045704399529
$ws.Range("J4").Value = 0.3691045704399529
$ws.Range("M4").Value = 0.07218766666666666
$ws.Range("N4").Value = 0.216563
$ws.Range("Q4").Value = 0.2054221089654445
$ws.Range("R4").Value = 1.848798980689
$ws.Range("S4").Value = 0.3691045704399529
$ws.Range("T4").Value = 0.3691045704399529

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.11543
$ws.Range("H5").Value = 0.34629
$ws.Range("I5").Value = 0.01497214206175765
$ws.Range("J5").Value = 0.01497214206175765
$ws.Range("M5").Value = 0.07218766666666666
$ws.Range("N5").Value = 0.216563
$ws.Range("Q5").Value = 0.008332622363333332
$ws.Range("R5").Value = 0.07499360127
$ws.Range("S5").Value = 0.01497214206175765
$ws.Range("T5").Value = 0.01497214206175765

# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 3.685254666666667
$ws.Range("H6").Value = 11.055764
$ws.Range("I6").Value = 0.4780053400596784
$ws.Range("J6").Value = 0.4780053400596784
$ws.Range("M6").Value = 0.07218766666666666
$ws.Range("N6").Value = 0.216563
$ws.Range("Q6").Value = 0.03626364278788889
$ws.Range("R6").Value = 0.326372785091
$ws.Range("S6").Value = 0.4780053400596784
$ws.Range("T6").Value = 0.4780053400596784
